$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new data row (row 9) with the 2024-01-02 22:58 resale numbers.
$ws.Range("A9").Value = "2024-01-02"
$ws.Range("B9").Value = "22:58:38"
$ws.Range("C9").Value = "Tuesday"
$ws.Range("D9").Value = "00"

$ws.Range("E9").Value = 140228
$ws.Range("F9").Value = 142835
$ws.Range("G9").Value = 171791
$ws.Range("H9").Value = 145980
$ws.Range("I9").Value = -1
$ws.Range("J9").Value = 117165
$ws.Range("K9").Value = 223709
$ws.Range("L9").Value = 248042
$ws.Range("M9").Value = 183700
$ws.Range("N9").Value = 109776
$ws.Range("O9").Value = 39778
$ws.Range("P9").Value = 30776
$ws.Range("Q9").Value = 71955
$ws.Range("R9").Value = -1
$ws.Range("S9").Value = 41683
$ws.Range("T9").Value = -1
